# Update the "Metadata" sheet (Property/Value table) with the new
# publication details (URL, Version, Date, Publisher).
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/participating-plan"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# The "Elements" sheet's Fixed Value cell for Extension.url mirrors the
# Metadata URL, so it picks up the new value automatically since it is the
# exact same shared text -- but to be safe against any divergence, set it
# explicitly too.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/participating-plan"

# The root "Extension" row (row 2) incorrectly repeated the
# ele-1/ext-1 constraint text in its Constraint(s) column; that text
# belongs only on the more specific rows (e.g. Extension.extension, row 4).
# Clear it here.
$elements.Range("AI2").Value = ""
